# Updated cryptos list with latest price and volume(1h) data
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '28.477.85'
$ws.Range('E2').Value = '  +1.35%  '
$ws.Range('D3').Value = '1.872.54'
$ws.Range('E3').Value = '  +1.30%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.009'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '315.86'
$ws.Range('E5').Value = '  +0.95%  '
$ws.Range('E6').Value = '  -0.12%  '
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3899'
$ws.Range('E8').Value = '  +0.18%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08349'
$ws.Range('E9').Value = '  +1.51%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.103'
$ws.Range('E10').Value = '  -0.10%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '41.74'
$ws.Range('E11').Value = '  +0.52%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '6.220'
$ws.Range('E12').Value = '  +0.78%  '
$ws.Range('D13').Value = '1.872.56'
$ws.Range('E13').Value = '  +0.93%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.39'
$ws.Range('E14').Value = '  +1.39%  '
$ws.Range('E15').Value = '  +1.58%  '
$ws.Range('E16').Value = '  -0.05%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001103'
$ws.Range('E17').Value = '  +0.84%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '91.15'
$ws.Range('E18').Value = '  +0.34%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06721'
$ws.Range('E19').Value = '  +0.76%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.67'
$ws.Range('E20').Value = '  +1.34%  '
$ws.Range('E21').Value = '  -0.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.902'
$ws.Range('E22').Value = '  +0.15%  '
$ws.Range('D23').Value = '28.508.43'
$ws.Range('E23').Value = '  +1.36%  '
$ws.Range('E24').Value = '  +1.03%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.228'
$ws.Range('E25').Value = '  -0.45%  '
$ws.Range('D26').Value = '2.087.50'
$ws.Range('E26').Value = '  +0.90%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '162.00'
$ws.Range('E27').Value = '  +1.60%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.60'
$ws.Range('E28').Value = '  +0.27%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.397'
$ws.Range('E29').Value = '  +1.48%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '125.84'
$ws.Range('E30').Value = '  +0.16%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.1042'
$ws.Range('E31').Value = '  +0.89%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.040'
$ws.Range('E32').Value = '  +2.34%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.748'
$ws.Range('E33').Value = '  -0.13%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.612'
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02451'
$ws.Range('E35').Value = '  +2.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06533'
$ws.Range('E36').Value = '  +2.21%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2155'
$ws.Range('E37').Value = '  -0.09%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.813'
$ws.Range('E38').Value = '  -2.45%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.018'
$ws.Range('E39').Value = '  +2.42%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.187'
$ws.Range('E40').Value = '  +1.34%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.240'
$ws.Range('E41').Value = '  -0.29%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.6376'
$ws.Range('E42').Value = '  +0.29%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.06'
$ws.Range('E43').Value = '  +0.42%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.006'
$ws.Range('E44').Value = '  -0.29%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5993'
$ws.Range('E45').Value = '  +0.88%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '13.02'
$ws.Range('E46').Value = '  +1.83%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.685'
$ws.Range('E47').Value = '  +0.38%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.005'
$ws.Range('E48').Value = '  +2.71%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.214'
$ws.Range('E49').Value = '  +1.83%  '
$ws.Range('E50').Value = '  +1.46%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.144'
$ws.Range('E51').Value = '  -10.19%  '
